# Edit derived from the target commit:
#  1. Slide 6's table switches to a different table style (tableStyleId).
#  2. The deck's theme colour scheme is swapped from the custom "Integral"
#     palette to the stock "Office Theme" palette (the two theme parts in
#     the package trade places; the colours are applied through the one
#     theme the PowerPoint object model exposes, since both theme parts
#     already share an identical font/format scheme and only differ in
#     their colour values and part-local `name` attributes).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{11E7D8CB-0BBE-4E57-B39B-B63CF18C471C}")

# --- 2. Theme colour scheme: Integral -> Office Theme --------------------
function Set-ThemeColor {
    param($scheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

Set-ThemeColor $colorScheme 1  "000000"   # dk1
Set-ThemeColor $colorScheme 2  "FFFFFF"   # lt1
Set-ThemeColor $colorScheme 3  "44546A"   # dk2
Set-ThemeColor $colorScheme 4  "E7E6E6"   # lt2
Set-ThemeColor $colorScheme 5  "5B9BD5"   # accent1
Set-ThemeColor $colorScheme 6  "ED7D31"   # accent2
Set-ThemeColor $colorScheme 7  "A5A5A5"   # accent3
Set-ThemeColor $colorScheme 8  "FFC000"   # accent4
Set-ThemeColor $colorScheme 9  "4472C4"   # accent5
Set-ThemeColor $colorScheme 10 "70AD47"   # accent6
Set-ThemeColor $colorScheme 11 "0563C1"   # hlink
Set-ThemeColor $colorScheme 12 "954F72"   # folHlink
